$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the Lydlus/password cell
#    (currently sitting between "Fm8AP" and "qpp") to right after
#    "Nath1234" in the Nath row, and merge the split
#    "Fm8AP" / "qpp" runs back into a single "Fm8APqpp" run.
# ---------------------------------------------------------------

# 1a) Insert the new (moved) bookmark right after "Nath1234".
#     A collapsed bookmark placed exactly at a paragraph-end
#     position can land in the wrong spot, so we temporarily type
#     a placeholder character after "Nath1234", anchor the new
#     bookmark just before it (a safe, non-boundary position), and
#     then remove the placeholder again.
$nathRng = $d.Content
$nathRng.Find.Execute("Nath1234", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nathRng.Collapse(0)
$nathRng.InsertAfter("X")

$safeSpot = $d.Range($nathRng.Start, $nathRng.Start)
$d.Bookmarks.Add("_GoBack", $safeSpot) | Out-Null

$placeholder = $d.Range($nathRng.Start, $nathRng.Start + 1)
$placeholder.Delete()

# 1b) Remove the old bookmark together with the text it sits in,
#     then retype the password so the two runs it used to split
#     ("Fm8AP" and "qpp") become a single "Fm8APqpp" run. Leaving a
#     single character behind before deleting the remainder avoids
#     emptying the paragraph completely (which would make the
#     insertion point of the retyped text land in the wrong cell).
#     "Fm8APqpp" (the visible, run-spanning text) shows up three
#     times in the table; only the Lydlus row (the last one) still
#     has the bookmark splitting it into two runs, so walk every
#     match and keep the last one found.
$lydlusRng = $d.Content
while ($lydlusRng.Find.Execute("Fm8APqpp", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastStart = $lydlusRng.Start
    $lastEnd = $lydlusRng.End
    $lydlusRng.Collapse(0)
    $lydlusRng.End = $d.Content.End
}
$passRng = $d.Range($lastStart + 1, $lastEnd)
$passRng.Delete()
$passRng.InsertAfter("m8APqpp")

# ---------------------------------------------------------------
# 2) Split "Yoan1234" into two runs: "Yoan123" stays as-is and a
#    new run containing "4" is appended with the same formatting.
#    A plain InsertAfter merges into the neighbouring run, so we
#    nudge a character-level property (Bold on, then back off) to
#    force the newly inserted text into its own run.
# ---------------------------------------------------------------
$yoanRng = $d.Content
$yoanRng.Find.Execute("Yoan123", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$yoanRng.Collapse(0)
$yoanRng.InsertAfter("4")
$yoanRng.Bold = 1
$yoanRng.Bold = 0

Write-Output "done"
